$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 1678.4166
$ws.Cells.Item(32, 9).Value = 1062.5
$ws.Cells.Item(32, 10).Value = 2294.3333
$ws.Cells.Item(32, 11).Value = 1062.5
$ws.Cells.Item(32, 12).Value = 2294.3333
$ws.Cells.Item(32, 13).Value = -736.5
$ws.Cells.Item(32, 14).Value = -2946.3333

$ws.Cells.Item(33, 8).Value = 576.1
$ws.Cells.Item(33, 9).Value = 612.3333
$ws.Cells.Item(33, 11).Value = 612.3333
$ws.Cells.Item(33, 13).Value = -383.3333

$ws.Cells.Item(61, 8).Value = 371.33334
$ws.Cells.Item(61, 9).Value = 371.33334
$ws.Cells.Item(61, 11).Value = 1114.00002
$ws.Cells.Item(61, 13).Value = -942.0000199999999

$ws.Cells.Item(76, 8).Value = 5874.5
$ws.Cells.Item(76, 9).Value = 5874.5
$ws.Cells.Item(76, 11).Value = 5874.5
$ws.Cells.Item(76, 13).Value = -5559.5

$ws.Cells.Item(79, 8).Value = 5874.5
$ws.Cells.Item(79, 9).Value = 5874.5
$ws.Cells.Item(79, 11).Value = 5874.5
$ws.Cells.Item(79, 13).Value = -4782.5

$ws.Cells.Item(101, 8).Value = 36500936
$ws.Cells.Item(101, 9).Value = 1001199
$ws.Cells.Item(101, 10).Value = 125250270
$ws.Cells.Item(101, 11).Value = 3003597
$ws.Cells.Item(101, 12).Value = 375750810
$ws.Cells.Item(101, 13).Value = -3001975
$ws.Cells.Item(101, 14).Value = -375754054

$ws.Cells.Item(113, 8).Value = 95317.63
$ws.Cells.Item(113, 9).Value = 500250
$ws.Cells.Item(113, 10).Value = 5332.6665
$ws.Cells.Item(113, 11).Value = 500250
$ws.Cells.Item(113, 12).Value = 5332.6665
$ws.Cells.Item(113, 13).Value = -496996
$ws.Cells.Item(113, 14).Value = -11840.6665

$ws.Cells.Item(127, 8).Value = 781
$ws.Cells.Item(127, 9).Value = 747.3333
$ws.Cells.Item(127, 11).Value = 2241.9999
$ws.Cells.Item(127, 13).Value = 2718.0001

$ws.Cells.Item(132, 8).Value = 5126.5
$ws.Cells.Item(132, 9).Value = 2476.4546
$ws.Cells.Item(132, 11).Value = 7429.3638
$ws.Cells.Item(132, 13).Value = -4899.3638

$ws.Cells.Item(138, 8).Value = 2855
$ws.Cells.Item(138, 9).Value = 3026.7273
$ws.Cells.Item(138, 10).Value = 2645.111
$ws.Cells.Item(138, 11).Value = 9080.1819
$ws.Cells.Item(138, 12).Value = 7935.333
$ws.Cells.Item(138, 13).Value = -3940.1819
$ws.Cells.Item(138, 14).Value = -18215.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2201
$ws.Cells.Item(61, 9).Value = 2201
$ws.Cells.Item(61, 11).Value = 2201
$ws.Cells.Item(61, 13).Value = -1989

$ws.Cells.Item(74, 8).Value = 4341.8667
$ws.Cells.Item(74, 9).Value = 4341.8667
$ws.Cells.Item(74, 11).Value = 4341.8667
$ws.Cells.Item(74, 13).Value = -3467.8667

$ws.Cells.Item(77, 8).Value = 4341.8667
$ws.Cells.Item(77, 9).Value = 4341.8667
$ws.Cells.Item(77, 11).Value = 21709.3335
$ws.Cells.Item(77, 13).Value = -17341.3335

$ws.Cells.Item(122, 8).Value = 3333.625
$ws.Cells.Item(122, 9).Value = 3395.8572
$ws.Cells.Item(122, 10).Value = 2898
$ws.Cells.Item(122, 11).Value = 10187.5716
$ws.Cells.Item(122, 12).Value = 8694
$ws.Cells.Item(122, 13).Value = -7737.571599999999
$ws.Cells.Item(122, 14).Value = -13594

$ws.Cells.Item(136, 8).Value = 2201
$ws.Cells.Item(136, 9).Value = 2201
$ws.Cells.Item(136, 11).Value = 6603
$ws.Cells.Item(136, 13).Value = -4053

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 4315.3
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 14).ClearContents()

$ws.Cells.Item(89, 8).Value = 4315.3
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 14).ClearContents()

$ws.Cells.Item(134, 8).Value = 116599.08
$ws.Cells.Item(134, 9).Value = 137298.95
$ws.Cells.Item(134, 11).Value = 411896.85
$ws.Cells.Item(134, 13).Value = -409361.85

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 62868.316
$ws.Cells.Item(31, 9).Value = 86931
$ws.Cells.Item(31, 11).Value = 86931
$ws.Cells.Item(31, 13).Value = -86636

$ws.Cells.Item(34, 8).Value = 62868.316
$ws.Cells.Item(34, 9).Value = 86931
$ws.Cells.Item(34, 11).Value = 86931
$ws.Cells.Item(34, 13).Value = -86729

$ws.Cells.Item(58, 8).Value = 1179
$ws.Cells.Item(58, 9).Value = 1252.6154
$ws.Cells.Item(58, 11).Value = 1252.6154
$ws.Cells.Item(58, 13).Value = -1049.6154

$ws.Cells.Item(136, 8).Value = 1179
$ws.Cells.Item(136, 9).Value = 1252.6154
$ws.Cells.Item(136, 11).Value = 3757.8462
$ws.Cells.Item(136, 13).Value = -1207.8462

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(64, 8).Value = 3470.4285
$ws.Cells.Item(64, 9).Value = 2764.3333
$ws.Cells.Item(64, 11).Value = 8292.999899999999
$ws.Cells.Item(64, 13).Value = -8022.999899999999

$ws.Cells.Item(67, 8).Value = 3470.4285
$ws.Cells.Item(67, 9).Value = 2764.3333
$ws.Cells.Item(67, 11).Value = 8292.999899999999
$ws.Cells.Item(67, 13).Value = -7356.999899999999

$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 14).ClearContents()

$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 14).ClearContents()

$ws.Cells.Item(88, 8).Value = 9195.777
$ws.Cells.Item(88, 10).Value = 9195.777
$ws.Cells.Item(88, 12).Value = 27587.331
$ws.Cells.Item(88, 14).Value = -28443.331

$ws.Cells.Item(91, 8).Value = 9195.777
$ws.Cells.Item(91, 10).Value = 9195.777
$ws.Cells.Item(91, 12).Value = 27587.331
$ws.Cells.Item(91, 14).Value = -30551.331

$ws.Cells.Item(94, 8).Value = 5143.75
$ws.Cells.Item(94, 10).Value = 5905.4
$ws.Cells.Item(94, 12).Value = 17716.2
$ws.Cells.Item(94, 14).Value = -19068.2

$ws.Cells.Item(101, 8).Value = 0
$ws.Cells.Item(101, 9).Value = 0
$ws.Cells.Item(101, 10).Value = 0
$ws.Cells.Item(101, 11).Value = 0
$ws.Cells.Item(101, 12).Value = 0
$ws.Cells.Item(101, 13).ClearContents()
$ws.Cells.Item(101, 14).ClearContents()

$ws.Cells.Item(104, 8).Value = 2242.5715
$ws.Cells.Item(104, 10).Value = 2166.3333
$ws.Cells.Item(104, 12).Value = 6498.999899999999
$ws.Cells.Item(104, 14).Value = -11740.9999

$ws.Cells.Item(114, 8).Value = 50000748
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 12).Value = 0
$ws.Cells.Item(114, 14).ClearContents()

$ws.Cells.Item(131, 8).Value = 1952.6875
$ws.Cells.Item(131, 9).Value = 1748.8334
$ws.Cells.Item(131, 11).Value = 5246.5002
$ws.Cells.Item(131, 13).Value = -206.5002000000004

$ws.Cells.Item(132, 8).Value = 2446.6667
$ws.Cells.Item(132, 9).Value = 1899.875
$ws.Cells.Item(132, 11).Value = 17098.875
$ws.Cells.Item(132, 13).Value = -14568.875

$ws.Cells.Item(137, 8).Value = 3354.6924
$ws.Cells.Item(137, 10).Value = 3741.5715
$ws.Cells.Item(137, 12).Value = 11224.7145
$ws.Cells.Item(137, 14).Value = -21424.7145

$ws.Cells.Item(138, 8).Value = 2756
$ws.Cells.Item(138, 9).Value = 1519.1666
$ws.Cells.Item(138, 10).Value = 6466.5
$ws.Cells.Item(138, 11).Value = 4557.4998
$ws.Cells.Item(138, 12).Value = 19399.5
$ws.Cells.Item(138, 13).Value = 582.5002000000004
$ws.Cells.Item(138, 14).Value = -29679.5

$ws.Cells.Item(141, 8).Value = 3869.2222
$ws.Cells.Item(141, 9).Value = 3046.2856
$ws.Cells.Item(141, 10).Value = 6749.5
$ws.Cells.Item(141, 11).Value = 9138.856800000001
$ws.Cells.Item(141, 12).Value = 20248.5
$ws.Cells.Item(141, 13).Value = -3958.856800000001
$ws.Cells.Item(141, 14).Value = -30608.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 7426.846
$ws.Cells.Item(132, 9).Value = 6002.9473
$ws.Cells.Item(132, 10).Value = 11291.714
$ws.Cells.Item(132, 11).Value = 18008.8419
$ws.Cells.Item(132, 12).Value = 33875.142
$ws.Cells.Item(132, 13).Value = -15478.8419
$ws.Cells.Item(132, 14).Value = -38935.142

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1001
$ws.Cells.Item(22, 9).Value = 1001
$ws.Cells.Item(22, 11).Value = 1001
$ws.Cells.Item(22, 13).Value = -706

$ws.Cells.Item(27, 8).Value = 1001
$ws.Cells.Item(27, 9).Value = 1001
$ws.Cells.Item(27, 11).Value = 1001
$ws.Cells.Item(27, 13).Value = -894

$ws.Cells.Item(46, 8).Value = 4669
$ws.Cells.Item(46, 9).Value = 4448.3335
$ws.Cells.Item(46, 11).Value = 4448.3335
$ws.Cells.Item(46, 13).Value = -4260.3335

$ws.Cells.Item(55, 8).Value = 371.85715
$ws.Cells.Item(55, 9).Value = 410.2
$ws.Cells.Item(55, 10).Value = 276
$ws.Cells.Item(55, 11).Value = 410.2
$ws.Cells.Item(55, 12).Value = 276
$ws.Cells.Item(55, 13).Value = -237.2
$ws.Cells.Item(55, 14).Value = -622

$ws.Cells.Item(132, 8).Value = 2949.2856
$ws.Cells.Item(132, 9).Value = 2588.4375
$ws.Cells.Item(132, 11).Value = 7765.3125
$ws.Cells.Item(132, 13).Value = -5235.3125

$ws.Cells.Item(136, 8).Value = 4941.6284
$ws.Cells.Item(136, 9).Value = 4684.4644
$ws.Cells.Item(136, 10).Value = 5970.2856
$ws.Cells.Item(136, 11).Value = 14053.3932
$ws.Cells.Item(136, 12).Value = 17910.8568
$ws.Cells.Item(136, 13).Value = -11503.3932
$ws.Cells.Item(136, 14).Value = -23010.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 50000
$ws.Cells.Item(5, 10).Value = 50000
$ws.Cells.Item(5, 12).Value = 50000
$ws.Cells.Item(5, 14).Value = -50224

$ws.Cells.Item(70, 8).Value = 20067.777
$ws.Cells.Item(70, 10).Value = 20067.777
$ws.Cells.Item(70, 12).Value = 20067.777
$ws.Cells.Item(70, 14).Value = -20697.777

$ws.Cells.Item(73, 8).Value = 20067.777
$ws.Cells.Item(73, 10).Value = 20067.777
$ws.Cells.Item(73, 12).Value = 20067.777
$ws.Cells.Item(73, 14).Value = -22251.777

$ws.Cells.Item(107, 8).Value = 16668027
$ws.Cells.Item(107, 10).Value = 50001310
$ws.Cells.Item(107, 12).Value = 150003930
$ws.Cells.Item(107, 14).Value = -150007770

$ws.Cells.Item(132, 8).Value = 14170.818
$ws.Cells.Item(132, 9).Value = 20697
$ws.Cells.Item(132, 11).Value = 62091
$ws.Cells.Item(132, 13).Value = -59561

$ws.Cells.Item(133, 8).Value = 65053.4
$ws.Cells.Item(133, 10).Value = 65053.4
$ws.Cells.Item(133, 12).Value = 65053.4
$ws.Cells.Item(133, 14).Value = -75173.39999999999
